# Updated cryptos list on Mon Jun 19 23:38:47 UTC 2023 with GitHub Actions
#
# The "cryptos" sheet lists the top coins in rows 2-51 (row 1 is the header):
#   A = rank index (unchanged), B = Coin, C = Link, D = Price, E = Volume(1h)
# This refresh pulls a newer price/volume snapshot. A few coins also swapped
# rank position (e.g. Toncoin/BitcoinCash swap places at rows 28-29, and a
# new coin - NEARProtocol - enters at the bottom row 51, displacing Frax and
# shifting rows 34-51 up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for columns B (Coin), C (Link), D (Price), E (Volume(1h))
# for every data row (2-51), taken from the refreshed feed.
$rows = @(
    @{ Row=2; Coin='Bitcoin'; Link='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; Price='26.721.90'; Volume='  +1.32%  ' },
    @{ Row=3; Coin='Ethereum'; Link='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; Price='1.732.86'; Volume='  +0.73%  ' },
    @{ Row=4; Coin='TetherUSD'; Link='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; Price='0.9980'; Volume='  -0.27%  ' },
    @{ Row=5; Coin='BNB'; Link='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; Price='242.36'; Volume='  -0.65%  ' },
    @{ Row=6; Coin='USDC'; Link='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; Price='0.9983'; Volume='  -0.26%  ' },
    @{ Row=7; Coin='XRP'; Link='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; Price='0.4928'; Volume='  +1.06%  ' },
    @{ Row=8; Coin='Cardano'; Link='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; Price='0.2626'; Volume='  +0.87%  ' },
    @{ Row=9; Coin='Dogecoin'; Link='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Price='0.06223'; Volume='  +0.40%  ' },
    @{ Row=10; Coin='WrappedEther'; Link='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Price='1.727.10'; Volume='  +0.32%  ' },
    @{ Row=11; Coin='Solana'; Link='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Price='15.93'; Volume='  +3.39%  ' },
    @{ Row=12; Coin='TRON'; Link='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Price='0.06980'; Volume='  -0.48%  ' },
    @{ Row=13; Coin='Polygon'; Link='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Price='0.6119'; Volume='  +2.81%  ' },
    @{ Row=14; Coin='Polkadot'; Link='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Price='4.506'; Volume='  -0.47%  ' },
    @{ Row=15; Coin='Litecoin'; Link='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Price='77.20'; Volume='  +0.04%  ' },
    @{ Row=16; Coin='Dai'; Link='https://coinranking.com/coin/MoTuySvg7+dai-dai'; Price='0.9982'; Volume='  -0.26%  ' },
    @{ Row=17; Coin='WrappedBTC'; Link='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Price='26.524.50'; Volume='  +0.51%  ' },
    @{ Row=18; Coin='BinanceUSD'; Link='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price='0.9981'; Volume='  -0.27%  ' },
    @{ Row=19; Coin='ShibaInu'; Link='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price='0.000007200'; Volume='  -0.11%  ' },
    @{ Row=20; Coin='Avalanche'; Link='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Price='11.42'; Volume='  +0.70%  ' },
    @{ Row=21; Coin='WrappedliquidstakedEther2.0'; Link='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Price='1.948.75'; Volume='  -0.03%  ' },
    @{ Row=22; Coin='Uniswap'; Link='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price='4.474'; Volume='  +0.01%  ' },
    @{ Row=23; Coin='Cosmos'; Link='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Price='8.566'; Volume='  +0.25%  ' },
    @{ Row=24; Coin='Chainlink'; Link='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Price='5.103'; Volume='  -1.16%  ' },
    @{ Row=25; Coin='Monero'; Link='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Price='138.36'; Volume='  +0.78%  ' },
    @{ Row=26; Coin='EthereumClassic'; Link='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Price='15.37'; Volume='  +0.74%  ' },
    @{ Row=27; Coin='LidoDAOToken'; Link='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price='1.771'; Volume='  +3.58%  ' },
    @{ Row=28; Coin='BitcoinCash'; Link='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Price='106.69'; Volume='  -0.89%  ' },
    @{ Row=29; Coin='Toncoin'; Link='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price='1.384'; Volume='  -2.39%  ' },
    @{ Row=30; Coin='InternetComputer(DFINITY)'; Link='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price='3.937'; Volume='  -0.40%  ' },
    @{ Row=31; Coin='Stellar'; Link='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price='0.07991'; Volume='  +0.61%  ' },
    @{ Row=32; Coin='Filecoin'; Link='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Price='3.666'; Volume='  -0.16%  ' },
    @{ Row=33; Coin='Hedera'; Link='https://coinranking.com/coin/jad286TjB+hedera-hbar'; Price='0.04487'; Volume='  -0.79%  ' },
    @{ Row=34; Coin='HuobiToken'; Link='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Price='2.608'; Volume='  -0.26%  ' },
    @{ Row=35; Coin='ARBITRUM'; Link='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Price='1.004'; Volume='  +1.15%  ' },
    @{ Row=36; Coin='ImmutableX'; Link='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price='0.6232'; Volume='  +0.46%  ' },
    @{ Row=37; Coin='TrustWalletToken'; Link='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Price='0.9412'; Volume='  +3.84%  ' },
    @{ Row=38; Coin='RenderToken'; Link='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Price='2.052'; Volume='  +3.90%  ' },
    @{ Row=39; Coin='MXToken'; Link='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Price='2.420'; Volume='  +0.98%  ' },
    @{ Row=40; Coin='PaxDollar'; Link='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; Price='0.9993'; Volume='  -0.07%  ' },
    @{ Row=41; Coin='VeChain'; Link='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Price='0.01512'; Volume='  +1.80%  ' },
    @{ Row=42; Coin='FraxShare'; Link='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Price='5.575'; Volume='  +3.64%  ' },
    @{ Row=43; Coin='Quant'; Link='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Price='99.49'; Volume='  -0.79%  ' },
    @{ Row=44; Coin='TheSandbox'; Link='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Price='0.3861'; Volume='  +0.52%  ' },
    @{ Row=45; Coin='Aptos'; Link='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Price='6.917'; Volume='  +2.86%  ' },
    @{ Row=46; Coin='Algorand'; Link='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Price='0.1161'; Volume='  +1.24%  ' },
    @{ Row=47; Coin='Cronos'; Link='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Price='0.05379'; Volume='  +0.44%  ' },
    @{ Row=48; Coin='EnergySwap'; Link='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price='7.895'; Volume='  +2.79%  ' },
    @{ Row=49; Coin='Elrond'; Link='https://coinranking.com/coin/omwkOTglq+elrond-egld'; Price='30.31'; Volume='  +0.96%  ' },
    @{ Row=50; Coin='Aave'; Link='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; Price='51.74'; Volume='  +1.57%  ' },
    @{ Row=51; Coin='NEARProtocol'; Link='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Price='1.231'; Volume='  -0.74%  ' }
)

foreach ($row in $rows) {
    $coinCell   = $ws.Cells.Item($row.Row, 2)
    $linkCell   = $ws.Cells.Item($row.Row, 3)
    $priceCell  = $ws.Cells.Item($row.Row, 4)
    $volumeCell = $ws.Cells.Item($row.Row, 5)

    if ($coinCell.Value -ne $row.Coin) {
        $coinCell.Value = $row.Coin
    }
    if ($linkCell.Value -ne $row.Link) {
        $linkCell.Value = $row.Link
    }
    if ($priceCell.Value -ne $row.Price) {
        # Price is stored as plain text ("0.9980", "242.36", "26.721.90", ...).
        # Several of those parse as valid numbers, and Excel's cell parser
        # would silently coerce them (dropping trailing zeros / adding
        # floating point noise) unless the cell is explicitly text-formatted
        # first.
        if ($row.Price -match '^-?[0-9]+(\.[0-9]+)?$') {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $row.Price
    }
    if ($volumeCell.Value -ne $row.Volume) {
        $volumeCell.Value = $row.Volume
    }
}
